$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.441.88'
$ws.Range('E2').Value = '  +3.74%  '

$ws.Range('D3').Value = '2.290.37'
$ws.Range('E3').Value = '  +3.48%  '

$ws.Range('E4').Value = '  -0.20%  '

$ws.Range('D5').Value = '320.89'
$ws.Range('E5').Value = '  +1.80%  '

$ws.Range('D6').Value = '105.36'
$ws.Range('E6').Value = '  +7.31%  '

$ws.Range('D7').Value = '0.592'
$ws.Range('E7').Value = '  +1.82%  '

$ws.Range('E8').Value = '  -0.18%  '

$ws.Range('D9').Value = '0.576'
$ws.Range('E9').Value = '  +2.98%  '

$ws.Range('D10').Value = '39.21'
$ws.Range('E10').Value = '  +7.49%  '

$ws.Range('D11').Value = '0.0845'
$ws.Range('E11').Value = '  +2.62%  '

$ws.Range('D12').Value = '7.96'
$ws.Range('E12').Value = '  +2.72%  '

$ws.Range('E13').Value = '  +2.24%  '

$ws.Range('D14').Value = '0.887'
$ws.Range('E14').Value = '  +3.02%  '

$ws.Range('D15').Value = '2.634.09'
$ws.Range('E15').Value = '  +3.39%  '

$ws.Range('D16').Value = '14.71'
$ws.Range('E16').Value = '  +4.19%  '

$ws.Range('D17').Value = '2.282.37'
$ws.Range('E17').Value = '  +3.65%  '

$ws.Range('D18').Value = '44.315.77'
$ws.Range('E18').Value = '  +3.86%  '

$ws.Range('D19').Value = '14.32'
$ws.Range('E19').Value = '  -3.06%  '

$ws.Range('E20').Value = '  +4.47%  '

$ws.Range('D21').Value = '6.61'
$ws.Range('E21').Value = '  +3.72%  '

$ws.Range('D22').Value = '66.55'
$ws.Range('E22').Value = '  +2.37%  '

$ws.Range('D23').Value = '3.24'
$ws.Range('E23').Value = '  +2.73%  '

$ws.Range('D24').Value = '238.64'
$ws.Range('E24').Value = '  +0.91%  '

$ws.Range('D25').Value = '2.22'
$ws.Range('E25').Value = '  +4.57%  '

$ws.Range('E26').Value = '  -0.20%  '

$ws.Range('D27').Value = '10.34'
$ws.Range('E27').Value = '  +3.17%  '

$ws.Range('D28').Value = '39.35'
$ws.Range('E28').Value = '  +16.21%  '

$ws.Range('E29').Value = '  +0.27%  '

$ws.Range('D30').Value = '6.59'
$ws.Range('E30').Value = '  +5.08%  '

$ws.Range('D31').Value = '163.82'
$ws.Range('E31').Value = '  +5.36%  '

$ws.Range('D32').Value = '0.0890'
$ws.Range('E32').Value = '  +1.85%  '

$ws.Range('D33').Value = '20.64'
$ws.Range('E33').Value = '  +1.22%  '

$ws.Range('D34').Value = '2.73'
$ws.Range('E34').Value = '  -0.81%  '

$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = '3.34'
$ws.Range('E35').Value = '  +5.74%  '

$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '2.11'
$ws.Range('E36').Value = '  +6.35%  '

$ws.Range('D37').Value = '0.115'
$ws.Range('E37').Value = '  +13.30%  '

$ws.Range('E38').Value = '  -0.40%  '

$ws.Range('D39').Value = '4.55'
$ws.Range('E39').Value = '  +2.81%  '

$ws.Range('D40').Value = '3.98'
$ws.Range('E40').Value = '  +7.69%  '

$ws.Range('D41').Value = '15.70'
$ws.Range('E41').Value = '  +29.98%  '

$ws.Range('E42').Value = '  +1.11%  '

$ws.Range('E43').Value = '  -0.19%  '

$ws.Range('D44').Value = '1.775.50'
$ws.Range('E44').Value = '  -5.25%  '

$ws.Range('E45').Value = '  +1.35%  '

$ws.Range('D46').Value = '86.23'
$ws.Range('E46').Value = '  -3.25%  '

$ws.Range('D47').Value = '5.44'
$ws.Range('E47').Value = '  +0.56%  '

$ws.Range('E48').Value = '  +4.13%  '

$ws.Range('D49').Value = '75.86'
$ws.Range('E49').Value = '  +0.25%  '

$ws.Range('D50').Value = '60.00'
$ws.Range('E50').Value = '  -0.47%  '

$ws.Range('D51').Value = '105.24'
$ws.Range('E51').Value = '  +4.28%  '
